$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are stored as plain text in this sheet (e.g. "317", "1048126.48").
# Writing a numeric-looking string via .Value normally gets auto-coerced to
# a real number, which would change the cell's stored type. Force text
# interpretation with a "@" number format while writing, then restore the
# original "Normal" style so no stray formatting is left behind.
$cells = @{
    "C9"  = "320"
    "D9"  = "300"
    "E9"  = "1088667.48"
    "C11" = "551"
    "D11" = "491"
    "E11" = "4231204.00"
    "C12" = "275"
    "E12" = "2184224.49"
    "C34" = "898"
    "E34" = "7313713.66"
    "C51" = "1181"
    "E51" = "9030924.05"
    "C52" = "814"
    "E52" = "5619450.79"
    "C60" = "6749"
    "E60" = "31007350.12"
    "C65" = "60"
    "D65" = "60"
    "E65" = "425050.00"
}

foreach ($addr in $cells.Keys) {
    $range = $ws.Range($addr)
    $range.NumberFormat = "@"
    $range.Value = $cells[$addr]
    $range.Style = "Normal"
}
